$wb = $excel.ActiveWorkbook

# Insert the new "After Loop" sheet right before the DATA sheet. We copy an
# existing sheet (instead of Worksheets.Add()) because freshly-added blank
# sheets in this runtime pick up a different default row height than sheets
# that were loaded from the original file; copying preserves the workbook's
# original sheet formatting (defaultRowHeight etc.) and we then overwrite
# its content.
$dataSheetBefore = $wb.Worksheets.Item("DATA")
$wb.Worksheets.Item("C").Copy($dataSheetBefore)

# Re-fetch sheet references by name now that indices have shifted because
# of the freshly inserted sheet.
$newSheet = $wb.Worksheets.Item("C (2)")
$newSheet.Name = "After Loop"
$cSheet = $wb.Worksheets.Item("C")
$dataSheet = $wb.Worksheets.Item("DATA")

# Wipe out the content/formatting that came along with the copy.
$newSheet.Cells.Clear()

# Header-ish row (row 2)
$newSheet.Range("A2").Value = "{{#each items}}{{value}}"
$newSheet.Range("B2").Value = "{{/each}}"
$newSheet.Range("A2:D2").Font.Color = 0

# Data / formula row (row 4) - cross-sheet formula references pointing at
# cells before/at/after the DATA sheet's used range.
$newSheet.Range("A4").Value = "{{#each items}}{{value}}"
$newSheet.Range("B4").Formula = "=A4+DATA!A1"
$newSheet.Range("C4").Formula = "=A4+DATA!B2"
$newSheet.Range("D4").Formula = "=A4+DATA!C3"
$newSheet.Range("E4").Value = "{{/each}}"

# Update selections on the other sheets first...
$cSheet.Range("B2").Select()

$dataSheet.Range("A3").Select()

# ...then leave the new "After Loop" sheet as the active/selected tab, since
# it becomes the active tab in the target workbook.
$newSheet.Range("C4").Select()

$wb.Save()
